# The output_processor sheet had an extraneous "topic_name" column (B) that
# duplicated/short-circuited the fuller "topic_description" column. This
# edit removes that column entirely so the sheet goes from describing each
# topic with a short name to describing it with its full description,
# shifting meeting_id and ingestion_timestep left by one column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting the entire column B shifts topic_description (C), meeting_id (D)
# and ingestion_timestep (E) left by one column (to B, C, D respectively),
# and automatically updates the sheet dimension from A1:E9 to A1:D9.
$ws.Range("B1").EntireColumn.Delete()
